# Update the "GRADUATE INFO SESSION" flyer image URL to point at the new
# image host link (the old officers-page/flyer image was removed, so the
# link was swapped for a working one), and leave the selection sitting on
# the cell that was actually edited (A2) instead of wherever it happened
# to be before (J18).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "https://i.ibb.co/kKjmFRH/grad-info-session-flyer.png"

$ws.Range("A2").Select()
